# Adds three new rows to the "Translation" sheet's text table (Table3,
# B3:I799) describing a new "speed / KM-H" readout, mirroring the existing
# RpmText "<value>" + number + Default-label triples already present
# (e.g. rows 74-76, 77-79, 80-82, 83-85) directly below the current last
# row (93).
#
# Columns: B=Text ID, C=Typography Name, D=Alignment, E=GB (display text),
#          F=Direction

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# Row 94: value placeholder text, centered, RpmText typography
$ws.Cells.Item(94, 2).Value = "SingleUseId97"
$ws.Cells.Item(94, 3).Value = "RpmText"
$ws.Cells.Item(94, 4).Value = "Center"
$ws.Cells.Item(94, 5).Value = "<value>"
$ws.Cells.Item(94, 6).Value = "LTR"

# Row 95: numeric-looking sample value "740" - it must stay a TEXT cell
# (matching the sheet's convention for every other numeric-looking sample,
# e.g. rows 75/78/81/84 = "21", rows 91/93 = "123"). A plain .Value/.Formula
# assignment of a digit string gets auto-converted to a real Number by the
# host, so instead enter it as a formula producing the text "740", then
# collapse the formula to its cached value via copy/paste-values - this
# keeps the literal text without leaving any NumberFormat/quote-prefix
# style behind.
$ws.Cells.Item(95, 2).Value = "SingleUseId98"
$ws.Cells.Item(95, 3).Value = "RpmText"
$ws.Cells.Item(95, 4).Value = "Left"
$ws.Cells.Item(95, 5).Formula = '="740"'
$ws.Cells.Item(95, 5).Copy()
$ws.Cells.Item(95, 5).PasteSpecial(-4163)
$ws.Cells.Item(95, 6).Value = "LTR"

# Row 96: the unit label shown next to the value
$ws.Cells.Item(96, 2).Value = "SingleUseId99"
$ws.Cells.Item(96, 3).Value = "Default"
$ws.Cells.Item(96, 4).Value = "Left"
$ws.Cells.Item(96, 5).Value = "KM/H"
$ws.Cells.Item(96, 6).Value = "LTR"
